# Applies the cryptos.xlsx data refresh described by the commit:
# "Updated cryptos list on Fri Aug 11 14:09:46 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.437.95"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "1.850.13"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'240.71"
$ws.Range("E5").Value = "  -0.78%  "
$ws.Range("D6").Value = "'0.6322"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "'0.9997"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("D9").Value = "'0.2956"
$ws.Range("E9").Value = "  -1.53%  "
$ws.Range("D10").Value = "'24.59"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").Value = "'0.07698"
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").Value = "1.859.52"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").Value = "'4.997"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("D14").Value = "'0.6866"
$ws.Range("E14").Value = "  -1.42%  "
$ws.Range("D15").Value = "'0.00001008"
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("D16").Value = "'83.20"
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("D17").Value = "2.107.55"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "'6.141"
$ws.Range("E18").Value = "  -2.26%  "
$ws.Range("D19").Value = "29.463.29"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("D20").Value = "'228.86"
$ws.Range("E20").Value = "  -2.89%  "
$ws.Range("D21").Value = "'12.51"
$ws.Range("E21").Value = "  -0.92%  "
$ws.Range("D22").Value = "'0.9996"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "'7.537"
$ws.Range("E23").Value = "  -1.94%  "
$ws.Range("D25").Value = "'157.04"
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("D26").Value = "'0.1397"
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("D27").Value = "'8.386"
$ws.Range("E27").Value = "  -1.82%  "
$ws.Range("E28").Value = "  -0.72%  "
$ws.Range("D29").Value = "'1.469"
$ws.Range("E29").Value = "  -0.84%  "
$ws.Range("D30").Value = "'1.270"
$ws.Range("E30").Value = "  +0.62%  "
$ws.Range("D31").Value = "'0.05704"
$ws.Range("E31").Value = "  -2.10%  "
$ws.Range("D32").Value = "'4.127"
$ws.Range("E32").Value = "  -0.44%  "
$ws.Range("D33").Value = "'4.032"
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("D34").Value = "'1.843"
$ws.Range("E34").Value = "  -3.55%  "
$ws.Range("D35").Value = "'1.157"
$ws.Range("E35").Value = "  -1.53%  "
$ws.Range("D36").Value = "'0.7148"
$ws.Range("E36").Value = "  -1.26%  "
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("D38").Value = "1.249.70"
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("D40").Value = "'2.778"
$ws.Range("E40").Value = "  -1.47%  "
$ws.Range("D41").Value = "'0.9077"
$ws.Range("E41").Value = "  -0.38%  "
$ws.Range("D42").Value = "'6.186"
$ws.Range("E42").Value = "  -0.48%  "
$ws.Range("D43").Value = "'1.0000"
$ws.Range("D44").Value = "'101.82"
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("D45").Value = "'66.15"
$ws.Range("E45").Value = "  -4.51%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "'0.00000000119"
$ws.Range("E46").Value = "  +1.27%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "'7.107"
$ws.Range("E47").Value = "  -3.80%  "
$ws.Range("D48").Value = "'0.4030"
$ws.Range("E48").Value = "  -1.08%  "
$ws.Range("D49").Value = "'9.100"
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("D50").Value = "'1.686"
$ws.Range("E50").Value = "  -2.34%  "
$ws.Range("D51").Value = "'0.1121"
$ws.Range("E51").Value = "  -0.33%  "
